$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-17 (columns A-T)
$data = @(
    ,@("ECs", "Agrn", "Atp1a3", "ECs", 3, 1, 15.140316, 45.420948, 0.3229157245229468, 0.3229157245229468, 1, 0.3333333333333333, 0.024647, 0.073941, 0.002001060000785908, 0.002001060000785908, 0.373163368452, 3.358470316068, 0.0006461737399676699, 0.0006461737399676698)
    ,@("ECs", "Agrn", "Atp1a3", "FAPs", 3, 1, 15.140316, 45.420948, 0.3229157245229468, 0.3229157245229468, 3, 1, 0.4009073333333333, 1.202722, 0.03254917956567031, 0.03254917956567031, 6.069863713384, 54.628773420456, 0.01051064190207592, 0.01051064190207592)
    ,@("ECs", "Agrn", "Atp1a3", "M2", 3, 1, 15.140316, 45.420948, 0.3229157245229468, 0.3229157245229468, 3, 1, 11.681588, 35.044764, 0.9484139445961233, 0.9484139445961232, 176.862933701808, 1591.766403316272, 0.306257776066923, 0.306257776066923)
    ,@("ECs", "Agrn", "Atp1a3", "sCs", 3, 1, 15.140316, 45.420948, 0.3229157245229468, 0.3229157245229468, 3, 1, 0.2098296666666667, 0.6294890000000001, 0.01703581583742065, 0.01703581583742065, 3.176887459508, 28.591987135572, 0.005501132813980181, 0.00550113281398018)
    ,@("FAPs", "Agrn", "Atp1a3", "ECs", 3, 1, 18.94069966666667, 56.822099, 0.4039710766824948, 0.4039710766824948, 1, 0.3333333333333333, 0.024647, 0.073941, 0.002001060000785908, 0.002001060000785908, 0.4668314246843334, 4.201482822159001, 0.0008083703630237572, 0.000808370363023757)
    ,@("FAPs", "Agrn", "Atp1a3", "FAPs", 3, 1, 18.94069966666667, 56.822099, 0.4039710766824948, 0.4039710766824948, 3, 1, 0.4009073333333333, 1.202722, 0.03254917956567031, 0.03254917956567031, 7.593465394830889, 68.341188553478, 0.0131489271142757, 0.01314892711427569)
    ,@("FAPs", "Agrn", "Atp1a3", "M2", 3, 1, 18.94069966666667, 56.822099, 0.4039710766824948, 0.4039710766824948, 3, 1, 11.681588, 35.044764, 0.9484139445961233, 0.9484139445961232, 221.2574499377373, 1991.317049439636, 0.3831318023391879, 0.3831318023391879)
    ,@("FAPs", "Agrn", "Atp1a3", "sCs", 3, 1, 18.94069966666667, 56.822099, 0.4039710766824948, 0.4039710766824948, 3, 1, 0.2098296666666667, 0.6294890000000001, 0.01703581583742065, 0.01703581583742065, 3.974320697490112, 35.76888627741101, 0.006881976866007518, 0.006881976866007517)
    ,@("M2", "Agrn", "Atp1a3", "ECs", 3, 1, 3.221232, 9.663696, 0.06870308817441464, 0.06870308817441464, 1, 0.3333333333333333, 0.024647, 0.073941, 0.002001060000785908, 0.002001060000785908, 0.07939370510400001, 0.7145433459360001, 0.0001374790016762885, 0.0001374790016762885)
    ,@("M2", "Agrn", "Atp1a3", "FAPs", 3, 1, 3.221232, 9.663696, 0.06870308817441464, 0.06870308817441464, 3, 1, 0.4009073333333333, 1.202722, 0.03254917956567031, 0.03254917956567031, 1.291415531168, 11.622739780512, 0.002236229153705103, 0.002236229153705103)
    ,@("M2", "Agrn", "Atp1a3", "M2", 3, 1, 3.221232, 9.663696, 0.06870308817441464, 0.06870308817441464, 3, 1, 11.681588, 35.044764, 0.9484139445961233, 0.9484139445961232, 37.629105076416, 338.661945687744, 0.06515896686143187, 0.06515896686143186)
    ,@("M2", "Agrn", "Atp1a3", "sCs", 3, 1, 3.221232, 9.663696, 0.06870308817441464, 0.06870308817441464, 3, 1, 0.2098296666666667, 0.6294890000000001, 0.01703581583742065, 0.01703581583742065, 0.6759100368160001, 6.083190331344, 0.001170413157601401, 0.0011704131576014)
    ,@("sCs", "Agrn", "Atp1a3", "ECs", 3, 1, 9.584029000000001, 28.752087, 0.2044101106201438, 0.2044101106201438, 1, 0.3333333333333333, 0.024647, 0.073941, 0.002001060000785908, 0.002001060000785908, 0.236217562763, 2.125958064867, 0.0004090368961181925, 0.0004090368961181925)
    ,@("sCs", "Agrn", "Atp1a3", "FAPs", 3, 1, 9.584029000000001, 28.752087, 0.2044101106201438, 0.2044101106201438, 3, 1, 0.4009073333333333, 1.202722, 0.03254917956567031, 0.03254917956567031, 3.842307508979334, 34.580767580814, 0.006653381395613593, 0.006653381395613592)
    ,@("sCs", "Agrn", "Atp1a3", "M2", 3, 1, 9.584029000000001, 28.752087, 0.2044101106201438, 0.2044101106201438, 3, 1, 11.681588, 35.044764, 0.9484139445961233, 0.9484139445961232, 111.956678158052, 1007.610103422468, 0.1938653993285805, 0.1938653993285805)
    ,@("sCs", "Agrn", "Atp1a3", "sCs", 3, 1, 9.584029000000001, 28.752087, 0.2044101106201438, 0.2044101106201438, 3, 1, 0.2098296666666667, 0.6294890000000001, 0.01703581583742065, 0.01703581583742065, 2.011013610393667, 18.099122493543, 0.003482292999831554, 0.003482292999831553)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $colNum = $j + 1
        [void]($ws.Cells.Item($rowNum, $colNum).Value = $rowData[$j])
    }
}